# Mise a jour avec pv presque complet 09.12.2025
# Adds two new parameter rows ("Puissance du module", "Surface du module") and a
# "Reviens retribution" row to the "Photovoltaic panel" table, fills in the
# previously-empty "Rendement" value, resizes the table, and tidies up the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Photovoltaic panel")

# --- Insert two new rows right above "Rendement" (currently row 4) ---------
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Row 4: Puissance du module
$ws.Cells.Item(4, 1).Value = "Puissance du module "
$ws.Cells.Item(4, 2).Value = 0.42
$ws.Cells.Item(4, 3).Value = $ws.Cells.Item(2, 3).Value2

# Row 5: Surface du module (unit is m^2, with superscript "2")
$ws.Cells.Item(5, 2).Value = 2
$unitCell = $ws.Cells.Item(5, 3)
$unitCell.Value = "m2"
# register the superscript font in the style table
$unitCell.Font.Superscript = $true
$unitCell.Font.Superscript = $false
# apply superscript only to the trailing "2" (rich-text run)
$unitCell.Characters(2, 1).Font.Superscript = $true

# Row 6 ("Rendement") now gets a value it was previously missing
$ws.Cells.Item(6, 2).Value = 25

# --- Insert a new row above "Emission de CO2" (currently row 10) -----------
$ws.Rows.Item(10).Insert()
$ws.Cells.Item(10, 1).Value = "Reviens retribution"
$ws.Cells.Item(10, 2).Value = 0.06
$ws.Cells.Item(10, 3).Value = $ws.Cells.Item(9, 3).Value2

# Label for row 5 is filled in last (matches the shared-string order of the
# authored workbook)
$ws.Cells.Item(5, 1).Value = "Surface du module"

# --- Resize the table to cover the new rows ---------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C13"))

# --- Column A widens to fit the new, longer labels --------------------------
$ws.Columns.Item(1).ColumnWidth = 17.5

# --- Restore the last selected cell ------------------------------------------
$ws.Range("C34").Select()
